# trigger Led Entfernt (#20)
# Trigger Led entfernt und wie Progress LED Pins umgelegt
#
# The "Farbsensor.trigger" row is removed; the three "progressLED" rows move
# up to take the place of LEDrot/LEDgruen/trigger, and LEDrot/LEDgruen move
# down below them. The last row in the block (37) becomes blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the "style blocks" first, using untouched template cells ---
# C20 already carries the "s7" (explicit black font) style used by the
# Farbsensor.LEDrot / Farbsensor.LEDgruen rows.
# C4 carries the plain/default style used by the Farbsensor.progressLEDx rows.
$ws.Range("C20").Copy()
$ws.Range("C35:C36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C4").Copy()
$ws.Range("C32:C34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Now move the text content into its new rows ---
$ws.Range("C32").Value = "Farbsensor.progressLED1"
$ws.Range("C33").Value = "Farbsensor.progressLED2"
$ws.Range("C34").Value = "Farbsensor.progressLED3"
$ws.Range("C35").Value = "Farbsensor.LEDrot"
$ws.Range("C36").Value = "Farbsensor.LEDgruen"

# --- Row 37 (was Farbsensor.trigger's old "progressLED3" slot) goes blank ---
$ws.Range("C37").ClearContents()
$ws.Range("B37").ClearContents()
$ws.Range("B37").Interior.ThemeColor = 2

# --- Selection, to match the saved view state ---
$ws.Range("C33").Select()
